$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three customer id values (shared strings) in column A
$ws.Range("A14").Value = "cus_NPuAUX7DbheAJC"
$ws.Range("A15").Value = "cus_NPuAJsxbM30H9R"
$ws.Range("A16").Value = "cus_NPuAexM6zAUeKZ"

# Normalize the current selection to a single cell A16
$ws.Range("A16").Select()
